$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 3 (B3, C3, D3 changed)
$ws.Range("B3").Value = 0.98095239999999995
$ws.Range("C3").Value = 0.63333329999999999
$ws.Range("D3").Value = 1

# Add new values for row 4 (E4, F4, G4)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.911111
$ws.Range("G4").Value = 0.99047620000000003

# Add new labels for the clustering section (rows 6 and 7)
$ws.Range("E6").Value = "SVM C = 1"
$ws.Range("B6").Value = "Alpha = 0.75"
$ws.Range("B7").Value = "Lambda = 0.25"

# Update the selected cell to B8
$ws.Range("B8").Select()
